$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new lesson row (#34 - Insert sort) by copying the formatting
# of the row above it (while it is still blank in column E) before
# filling in the new values, so the new row inherits the same column
# styles (date format, wrap text, etc.) without picking up row 30's
# soon-to-be-added payment value.
$ws.Rows(30).Copy()
$ws.Rows(31).PasteSpecial()
$excel.CutCopyMode = $false

$ws.Range("A31").Value = 34
$ws.Range("B31").Value = 44572
$ws.Range("C31").Value = "Insert sort"
$ws.Range("D31").Value = "D:\Teaching\12.Algorithms\Tasks\Tasks"

# Fill in the payment column for the last few already-taught lessons
# (rows 27-30), which had been left blank until now.
$ws.Range("E27").Value = 800
$ws.Range("E28").Value = 800
$ws.Range("E29").Value = 800
$ws.Range("E30").Value = 800

# Match the author's final selection/scroll state.
$ws.Range("E31").Select()
